$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.222.76'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.839.17'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.81'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('E5').NumberFormat = 'General'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4676'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.31%  '
$ws.Range('E7').NumberFormat = 'General'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2709'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.34%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06276'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.48%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.838.89'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07409'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.12'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.928'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.96%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '83.73'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.95%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6201'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.43%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.143.11'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9999'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.00'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007298'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('E19').NumberFormat = 'General'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.36'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.80%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9996'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.886'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.05%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.848'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.206'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('E24').NumberFormat = 'General'
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.35'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.59%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.81'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.19%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.873'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('E27').NumberFormat = 'General'
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1040'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.370'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.078'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.63%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.798'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.86%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04814'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.139'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.38%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7092'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('E34').NumberFormat = 'General'
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.697'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01873'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.28%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.643'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('E37').NumberFormat = 'General'
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8922'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.921'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.08%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '104.46'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.002'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.542'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4013'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.51%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.041'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1193'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.88%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.83'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.509'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.57%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.77'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05510'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.356'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.06%  '
$ws.Range('E50').NumberFormat = 'General'
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3639'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.21%  '
$ws.Range('E51').NumberFormat = 'General'
